$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.718.93"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "1.885.93"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4723"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3972"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08064"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").Value = "1.884.49"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.974"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.211"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001041"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06592"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").Value = "27.731.97"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.522"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.309"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "2.102.06"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.601"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9676"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09537"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.477"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.623"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.311"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06127"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02256"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.188"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6010"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1898"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.79%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5700"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.404"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.943"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06829"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.45%  "
